$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while forcing text storage
# (so numeric-looking strings like "0.614" or "55.90" are not silently
# auto-converted to numbers / lose formatting), and then clear the
# temporary text-number-format back to the sheet's default style so no
# stray cell style (s="n") attribute is left behind.
function Set-TextValue {
    param($addr, $value)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "37.327.91"
$ws.Range("E2").Value = "  +0.32%  "
Set-TextValue "D3" "2.006.94"
$ws.Range("E3").Value = "  +0.12%  "
Set-TextValue "D5" "257.33"
$ws.Range("E5").Value = "  +4.33%  "
Set-TextValue "D6" "0.614"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  -0.02%  "
Set-TextValue "D8" "55.90"
$ws.Range("E8").Value = "  -6.39%  "
$ws.Range("E9").Value = "  -1.01%  "
Set-TextValue "D10" "0.0767"
$ws.Range("E10").Value = "  -5.46%  "
$ws.Range("E11").Value = "  -1.86%  "
Set-TextValue "D12" "2.304.12"
$ws.Range("E12").Value = "  +0.24%  "
Set-TextValue "D13" "14.20"
$ws.Range("E13").Value = "  -5.49%  "
Set-TextValue "D14" "21.08"
$ws.Range("E14").Value = "  -5.40%  "
Set-TextValue "D15" "0.799"
$ws.Range("E15").Value = "  -5.37%  "
Set-TextValue "D16" "5.22"
$ws.Range("E16").Value = "  -4.42%  "
Set-TextValue "D17" "2.008.19"
$ws.Range("E17").Value = "  +0.00%  "
Set-TextValue "D18" "37.109.62"
$ws.Range("E18").Value = "  +0.10%  "
Set-TextValue "D19" "69.53"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  -3.80%  "
$ws.Range("E21").Value = "  -1.62%  "
Set-TextValue "D22" "228.20"
$ws.Range("E22").Value = "  -0.83%  "
Set-TextValue "D23" "2.60"
$ws.Range("E23").Value = "  +4.53%  "
Set-TextValue "D24" "0.999"
$ws.Range("E24").Value = "  -0.18%  "
Set-TextValue "D25" "2.34"
$ws.Range("E25").Value = "  -0.05%  "
Set-TextValue "D26" "164.60"
$ws.Range("E26").Value = "  +0.17%  "
Set-TextValue "D27" "8.86"
$ws.Range("E27").Value = "  -6.23%  "
Set-TextValue "D28" "19.59"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -7.57%  "
$ws.Range("E30").Value = "  -3.24%  "
$ws.Range("E31").Value = "  -1.51%  "
Set-TextValue "D32" "4.62"
$ws.Range("E32").Value = "  -3.83%  "
Set-TextValue "D33" "0.0641"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("E35").Value = "  -2.49%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  -0.07%  "
Set-TextValue "D38" "3.35"
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").Value = "  +3.80%  "
Set-TextValue "D41" "1.20"
$ws.Range("E41").Value = "  +1.45%  "
Set-TextValue "D42" "0.0928"
$ws.Range("E42").Value = "  -5.54%  "
Set-TextValue "D43" "0.0211"
$ws.Range("E43").Value = "  -1.21%  "
Set-TextValue "D44" "1.399.29"
$ws.Range("E44").Value = "  +2.24%  "
Set-TextValue "D45" "89.40"
$ws.Range("E45").Value = "  -2.94%  "
Set-TextValue "D46" "15.67"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D48" "2.90"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D49" "7.01"
$ws.Range("E49").Value = "  -4.88%  "
Set-TextValue "D50" "2.196.18"
$ws.Range("E50").Value = "  +0.27%  "
Set-TextValue "D51" "1.93"
$ws.Range("E51").Value = "  -6.77%  "
